$d = $word.ActiveDocument

# The runtime's Range.Find.Execute(... Replace:=wdReplaceOne) does not
# reliably confine its match to the calling Range (it can hit the first
# occurrence anywhere in the document), which is unsafe here because
# several cells briefly share identical old/new text during this edit.
# So cell/paragraph text is rewritten directly by Range position instead,
# trimming the trailing cell-mark / paragraph-mark character off first.

# Update the date line at the top of the document
$dateRange = $d.Paragraphs.Item(1).Range
$d.Range($dateRange.Start, $dateRange.End - 1).Text = "2025-01-25 Saturday"

# Update each division-problem cell in the table
$t = $d.Tables.Item(1)

$c = $t.Cell(1, 1).Range
$d.Range($c.Start, $c.End - 1).Text = "94÷5="  # was "19÷4="
$c = $t.Cell(1, 2).Range
$d.Range($c.Start, $c.End - 1).Text = "60÷7="  # was "42÷2="
$c = $t.Cell(1, 3).Range
$d.Range($c.Start, $c.End - 1).Text = "58÷2="  # was "90÷8="
$c = $t.Cell(1, 4).Range
$d.Range($c.Start, $c.End - 1).Text = "37÷7="  # was "32÷9="
$c = $t.Cell(1, 5).Range
$d.Range($c.Start, $c.End - 1).Text = "16÷7="  # was "63÷4="

$c = $t.Cell(5, 1).Range
$d.Range($c.Start, $c.End - 1).Text = "70÷4="  # was "47÷5="
$c = $t.Cell(5, 2).Range
$d.Range($c.Start, $c.End - 1).Text = "56÷5="  # was "58÷2="
$c = $t.Cell(5, 3).Range
$d.Range($c.Start, $c.End - 1).Text = "41÷4="  # was "28÷2="
$c = $t.Cell(5, 4).Range
$d.Range($c.Start, $c.End - 1).Text = "13÷7="  # was "61÷5="
$c = $t.Cell(5, 5).Range
$d.Range($c.Start, $c.End - 1).Text = "13÷6="  # was "34÷8="

$c = $t.Cell(9, 1).Range
$d.Range($c.Start, $c.End - 1).Text = "17÷6="  # was "89÷2="
$c = $t.Cell(9, 2).Range
$d.Range($c.Start, $c.End - 1).Text = "83÷7="  # was "72÷4="
$c = $t.Cell(9, 3).Range
$d.Range($c.Start, $c.End - 1).Text = "50÷7="  # was "99÷4="
$c = $t.Cell(9, 4).Range
$d.Range($c.Start, $c.End - 1).Text = "58÷8="  # was "53÷8="
$c = $t.Cell(9, 5).Range
$d.Range($c.Start, $c.End - 1).Text = "46÷8="  # was "64÷3="

$c = $t.Cell(13, 1).Range
$d.Range($c.Start, $c.End - 1).Text = "25÷4="  # was "42÷8="
$c = $t.Cell(13, 2).Range
$d.Range($c.Start, $c.End - 1).Text = "76÷2="  # was "27÷6="
$c = $t.Cell(13, 3).Range
$d.Range($c.Start, $c.End - 1).Text = "45÷9="  # was "85÷8="
$c = $t.Cell(13, 4).Range
$d.Range($c.Start, $c.End - 1).Text = "60÷3="  # was "67÷5="
$c = $t.Cell(13, 5).Range
$d.Range($c.Start, $c.End - 1).Text = "24÷8="  # was "86÷4="

$c = $t.Cell(17, 1).Range
$d.Range($c.Start, $c.End - 1).Text = "83÷5="  # was "21÷3="
$c = $t.Cell(17, 2).Range
$d.Range($c.Start, $c.End - 1).Text = "25÷3="  # was "32÷5="
$c = $t.Cell(17, 3).Range
$d.Range($c.Start, $c.End - 1).Text = "39÷2="  # was "44÷6="
$c = $t.Cell(17, 4).Range
$d.Range($c.Start, $c.End - 1).Text = "63÷2="  # was "62÷6="
$c = $t.Cell(17, 5).Range
$d.Range($c.Start, $c.End - 1).Text = "31÷3="  # was "28÷5="

Write-Host "Replacements complete"
